$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change zoom scale (85% -> 115%)
$excel.ActiveWindow.Zoom = 115

# Data: Cilium_number, horizontal_length_in_pixels, z_lower, z_upper
$data = @(
    @(1, 10.6, 2, 9),
    @(2, 11.364000000000001, 3, 10),
    @(3, 13.276, 2, 8),
    @(4, 19.771999999999998, 3, 9),
    @(5, 15.7, 3, 10),
    @(6, 5.3360000000000003, 6, 12),
    @(7, 15.042, 4, 12),
    @(8, 11.241, 4, 10),
    @(9, 14.535, 3, 10),
    @(10, 16.986000000000001, 3, 11),
    @(11, 9.4309999999999992, 5, 11),
    @(12, 19.800999999999998, 3, 10),
    @(13, 19.648, 4, 11),
    @(14, 25.466000000000001, 3, 10),
    @(15, 18.992000000000001, 4, 11),
    @(16, 16.108000000000001, 5, 10)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
